$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestSuite1")

# ---------------------------------------------------------------------------
# 1. Sheet1 ("TestSuite1"): add two new test-case rows (5 and 6) and a
#    trailing formatted-only row (7), reusing the existing "row3" formatting
#    (plain wrap-text style) as a template so no stray styles are created.
# ---------------------------------------------------------------------------
$ws1.Range("A3:AD3").Copy()
$ws1.Range("A5:AD5").PasteSpecial(-4122)
$ws1.Range("A6:AD6").PasteSpecial(-4122)

$ws1.Range("A5").Value = 'TC04'
$ws1.Range("B5").Value = 'WalletRegistration'
$ws1.Range("C5").Value = 'P_WalletRegister'
$ws1.Range("D5").Value = 'Create new customer- All valid details'
$ws1.Range("E5").Value = 'POST'
$ws1.Range("F5").Value = '/wallet/v1/customers/create'
$ws1.Range("G5").Value = '{ 
"mobileNo":#new_mobNo, 
"email":#new_email_id,
"name_of_customer":"TestUser"
}'
$ws1.Range("H5").Value = '{
"status": "success",
"status_msg": "Customer created successfully",
"consumer_id": $$new_cust_id
}'

$ws1.Range("A6").Value = 'TC05'
$ws1.Range("B6").Value = 'WalletRegistration'
$ws1.Range("C6").Value = 'P_WalletRegister'
$ws1.Range("D6").Value = 'Create new customer- All valid details'
$ws1.Range("E6").Value = 'POST'
$ws1.Range("F6").Value = '/wallet/v1/customers/create'
$ws1.Range("G6").Value = '{ 
"mobileNo":"9833868977", 
"email":"qc9167916845@gmail.com",
"name_of_customer":"WLWTestTwo"
}'
$ws1.Range("H6").Value = '{
"status": "success",
"status_msg": "Customer created successfully",
"consumer_id": $$new_cust_id
}'

# Expected status-code column gets a right-aligned variant of the same style.
$ws1.Range("I5:I6").WrapText = $true
$ws1.Range("I5:I6").HorizontalAlignment = -4152
$ws1.Range("I5").Value = 201
$ws1.Range("I6").Value = 201

$ws1.Rows.Item(5).RowHeight = 84.75
$ws1.Rows.Item(6).RowHeight = 84.75

# Row 7: just a single formatted (blank) cell at G7, matching the template.
$ws1.Range("J4").Copy()
$ws1.Range("G7").PasteSpecial(-4122)

$ws1.Range("E2").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "TestData" sheet straight after "TestSuite1".
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Add($null, $ws1)
$wsData.Name = "TestData"

$ws1.Range("J4").Copy()
$wsData.Range("A1:G4").PasteSpecial(-4122)

$wsData.Range("A2:A3").WrapText = $true
$wsData.Range("A2:A3").HorizontalAlignment = -4152
$wsData.Range("C2:D3").WrapText = $true
$wsData.Range("C2:D3").HorizontalAlignment = -4152
$wsData.Range("G2:G4").WrapText = $true
$wsData.Range("G2:G4").HorizontalAlignment = -4152

$wsData.Range("A1").Value = '#new_mobNo'
$wsData.Range("B1").Value = 'new_cust_id'
$wsData.Range("C1").Value = 'old_mobNo'
$wsData.Range("D1").Value = 'old_cust_id'
$wsData.Range("E1").Value = '#new_email_id'
$wsData.Range("F1").Value = 'old_email_id'
$wsData.Range("G1").Value = 'cust_id_zero_bal'

$wsData.Range("A2").Value = 9810101010
$wsData.Range("C2").Value = 9833868977
$wsData.Range("D2").Value = 833012
$wsData.Range("E2").Value = 'abc@bc.com'
$wsData.Range("F2").Value = 'test@gmail.com'
$wsData.Range("G2").Value = 123

$wsData.Range("A3").Value = 9812231223
$wsData.Range("C3").Value = 7720077155
$wsData.Range("D3").Value = 143245
$wsData.Range("G3").Value = 34567

$wsData.Range("G4").Value = 1234567

$wsData.Rows.Item(2).RowHeight = 24.75

$wsData.Columns.Item(1).ColumnWidth = 13.59244791666667
$wsData.Columns.Item(2).ColumnWidth = 12.87760416666667
$wsData.Columns.Item(3).ColumnWidth = 15.59244791666667
$wsData.Columns.Item(4).ColumnWidth = 15.30729166666667
$wsData.Columns.Item(5).ColumnWidth = 17.87760416666667
$wsData.Columns.Item(6).ColumnWidth = 14.30729166666667
$wsData.Columns.Item(7).ColumnWidth = 16.59244791666667

$wsData.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. Restore TestSuite1 as the active sheet/tab (it stays selected after the
#    edit, per the original authoring session).
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Select()
